$d = $word.ActiveDocument

# ------------------------------------------------------------------------
# Edit 1: Heading "Education / Coursework" -> "Former Education / Coursework"
#         A new "Former " run is inserted right before "Education".
# ------------------------------------------------------------------------
$heading = $d.Content
$heading.Find.Execute("Education / Coursework") | Out-Null
$insertPoint = $heading.Duplicate
$insertPoint.Collapse(1)                 # wdCollapseStart
$insertPoint.InsertBefore("Former ")

# ------------------------------------------------------------------------
# Edit 2: Education bullet list is re-shuffled:
#   "B.S. In ECE | May 2020"                   gains the 16 trailing spaces
#   "M.S. In ECE | May 2021" + 16 spaces    -> "Dean's List Spring 2018"
#   "Dean's List Spring 2018 "               -> "Stuyvesant H.S. Class of 2016"
#
# Order matters: the "Dean's List Spring 2018 " replacement is done first,
# while that phrase is still unique in the document; only afterwards do we
# turn the old "M.S." line into the (now unique again) "Dean's List" text,
# so a single Find.Execute/Replace never matches more than one paragraph.
# ------------------------------------------------------------------------
$d.Content.Find.Execute("Dean" + [char]0x2019 + "s List Spring 2018 ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Stuyvesant H.S. Class of 2016", 2) | Out-Null

$d.Content.Find.Execute("M.S. In ECE | May 2021", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Dean" + [char]0x2019 + "s List Spring 2018", 2) | Out-Null

$bsHit = $d.Content
$bsHit.Find.Execute("B.S. In ECE | May 2020") | Out-Null
$bsParagraph = $bsHit.Paragraphs(1).Range
$bsParagraph.InsertAfter("                ")

# ------------------------------------------------------------------------
# Edit 3: The "_GoBack" bookmark moves from in front of "| AWS | Azure |"
#         to between "graduate students " and "fundamental ".
# ------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
  $d.Bookmarks("_GoBack").Delete()
}
$bookmarkAnchor = $d.Content
$bookmarkAnchor.Find.Execute("graduate students ") | Out-Null
$bookmarkPoint = $d.Range($bookmarkAnchor.End, $bookmarkAnchor.End)
$d.Bookmarks.Add("_GoBack", $bookmarkPoint)

Write-Output "Applied resume edits."
